# Generate Report for Handback
#
# This applies the "handback" localization report refresh:
#  - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#  - zh-cn sheet: populate Latest Target File / Latest Handback File (+ hyperlink)
#    and set the Latest Handback DateTime
#  - de-de sheet: same, with its own (later) handback datetime
#  - Widen a few columns so the new/longer values are readable

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

$newStatus = "Handed back: in sync with en-US"

# --- 1. Status text everywhere it appears -----------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- 2. zh-cn sheet: Latest Target File / Latest Handback File --------
$zhRow2Md  = "335d0462-a42d-4902-b4f9-61180468c70a.md"
$zhRow3Md  = "8e92d911-80b0-43a1-a734-d851eada6e9b.md"
$zhRow2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6c453d5f2d054f3e6b39f2e28d67f424d8dc655d/e2e/335d0462-a42d-4902-b4f9-61180468c70a.md"
$zhRow3Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6c453d5f2d054f3e6b39f2e28d67f424d8dc655d/e2e/8e92d911-80b0-43a1-a734-d851eada6e9b.md"

$wsZhCn.Range("J2").Value = "335d0462-a42d-4902-b4f9-61180468c70a.a12ba7382b1f953382826405c2b6c6904426ad0e.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "8e92d911-80b0-43a1-a734-d851eada6e9b.d4b366fdd5791641094334132cc0538a5d35861b.zh-cn.xlf"

$wsZhCn.Range("K2").Value = "2016-08-28 06:35:11"
$wsZhCn.Range("K3").Value = "2016-08-28 06:35:11"

# Rebuild hyperlinks on the zh-cn sheet in the final order: A2, I2, A3, I3
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhRow2Url, "", "", $zhRow2Md)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhRow2Url, "", "", $zhRow2Md)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $zhRow3Url, "", "", $zhRow3Md)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhRow3Url, "", "", $zhRow3Md)

# --- 3. de-de sheet: Latest Target File / Latest Handback File --------
$deRow2Md  = "335d0462-a42d-4902-b4f9-61180468c70a.md"
$deRow3Md  = "8e92d911-80b0-43a1-a734-d851eada6e9b.md"
$deRow2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6c453d5f2d054f3e6b39f2e28d67f424d8dc655d/e2e/335d0462-a42d-4902-b4f9-61180468c70a.md"
$deRow3Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6c453d5f2d054f3e6b39f2e28d67f424d8dc655d/e2e/8e92d911-80b0-43a1-a734-d851eada6e9b.md"

$wsDeDe.Range("J2").Value = "335d0462-a42d-4902-b4f9-61180468c70a.a12ba7382b1f953382826405c2b6c6904426ad0e.de-de.xlf"
$wsDeDe.Range("J3").Value = "8e92d911-80b0-43a1-a734-d851eada6e9b.d4b366fdd5791641094334132cc0538a5d35861b.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-08-28 06:35:17"
$wsDeDe.Range("K3").Value = "2016-08-28 06:35:17"

# Rebuild hyperlinks on the de-de sheet in the final order: A2, I2, A3, I3
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deRow2Url, "", "", $deRow2Md)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $deRow2Url, "", "", $deRow2Md)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $deRow3Url, "", "", $deRow3Md)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $deRow3Url, "", "", $deRow3Md)

# --- 4. Column widths ---------------------------------------------------
# ColumnWidth inputs are chosen so the engine's internal pixel-rounding
# lands on the widest achievable value closest to the target (status text
# got much longer; the new file-name/hyperlink columns need to show full
# GUID-based file names).
$wsOverview.Columns.Item(5).ColumnWidth = 29.1   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = 29.1   # F: de-de status

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.1   # C: Status
$wsZhCn.Columns.Item(9).ColumnWidth  = 39.1   # I: Latest Target File
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1   # J: Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.1   # C: Status
$wsDeDe.Columns.Item(9).ColumnWidth  = 39.1   # I: Latest Target File
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1   # J: Latest Handback File
